$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; everything currently at row 31+ shifts down one.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Range("A31").Value = 10
$ws.Range("B31").Value = "Vega Modelo de Temuco"
$ws.Range("C31").Value = "La Araucanía"
$ws.Range("D31").Value = 44629
$ws.Range("E31").Value = 9
$ws.Range("F31").Value = 100114002
$ws.Range("G31").Value = "Camote"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 30
$ws.Range("K31").Value = 18000
$ws.Range("L31").Value = 18000
$ws.Range("M31").Value = 18000
$ws.Range("N31").Value = "$/malla 20 kilos"
$ws.Range("O31").Value = "Perú"
$ws.Range("P31").Value = 900
$ws.Range("Q31").Value = 20
$ws.Range("R31").Value = "Hortaliza"
